$d = $word.ActiveDocument

# --- Step 1: remove the stray "_GoBack" bookmark currently sitting after
#     "Age significantly improves vocab" at the very end of the document ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: fix the "differnces" -> "differences" typo and drop the
#     now-stale spell-check proofErr markers that bracketed it, then plant
#     a fresh (collapsed) "_GoBack" bookmark right after the corrected word ---

# locate the word that needs fixing
$find = $d.Content
$find.Find.Execute("differnces", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPos = $find.Start

# resolve the *paragraph* that contains it (walking $d.Paragraphs keeps this
# robust to any future bookkeeping/indexing differences)
$paraCount = $d.Paragraphs.Count
$paraRange = $null
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range
    if ($targetPos -ge $candidate.Start -and $targetPos -lt $candidate.End) {
        $paraRange = $candidate
        break
    }
}

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="001D2CE2" w:rsidRDefault="001D2CE2" w:rsidP="00490FE0"><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="360"/><w:tab w:val="num" w:pos="1080"/></w:tabs><w:ind w:left="1080"/></w:pPr><w:r><w:t xml:space="preserve">Mean of </w:t></w:r><w:r><w:t>differences</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> mu(d) = mu(2) – mu(1) = 12 – 3 = </w:t></w:r><w:r w:rsidRPr="00667BDA"><w:rPr><w:b/></w:rPr><w:t>9</w:t></w:r></w:p>
'@

$paraRange.InsertXML($xml)
